# Implemented Remote Driver through driverFactory
# Updates the "TestData" sheet: fills in previously blank UserName/Password
# and HamburgerMenuOption/HamburgerSubMenu/ExpectedPageTitle cells with
# "Test"/"test" placeholder values, and moves the active selection.
# A leading apostrophe forces these into the same quote-prefixed text style
# ("s=2" / wrap-text with quotePrefix) the surrounding blank cells already
# use, instead of Excel re-picking a plain style when the value is assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 2 (login / No / Chrome / Admin / admin123 ...)
$ws.Range("F2").Value = "'Test"
$ws.Range("G2").Value = "'test"
$ws.Range("H2").Value = "'test"

# Row 3 (login / No / Firefox / Admin123 / admin123 ...)
$ws.Range("F3").Value = "'Test"
$ws.Range("G3").Value = "'test"
$ws.Range("H3").Value = "'test"

# Row 4 (verifyAmazonHamburgerSubMenuPageTitle / Yes / Chrome ...)
$ws.Range("D4").Value = "'Test"
$ws.Range("E4").Value = "'test"

# Update the active selection to match the authored workbook
$ws.Range("G7").Select()
